# Auto-generated Word COM-interop script implementing the target diff.
# Strategy: locate each target paragraph by its distinctive anchor text,
# delete the paragraph (including its trailing paragraph mark), then
# InsertXML the precise replacement run/paragraph structure extracted
# programmatically from the unified diff. InsertXML preserves exact run
# boundaries (even among identically-formatted adjacent runs), unlike
# Find/Replace or Range.Text which normalizes/merges same-format runs.

$d = $word.ActiveDocument

function Get-ParagraphByStart([string]$anchor) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        if ($t.StartsWith($anchor)) {
            return $p
        }
    }
    return $null
}

# --- Edit 1: "And one more time my karma appeared..." paragraph ---
$target1 = Get-ParagraphByStart("And one more time my karma appeared")
if ($null -eq $target1) {
    throw "Could not locate paragraph 1 (And one more time...)"
}
$r1 = $target1.Range
$delRange1 = $d.Range($r1.Start, $r1.End + 1)
$delRange1.Delete()

$insPoint1 = $d.Range($r1.Start, $r1.Start)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">And one more time, my karma appeared. Ouh, right, I said that I don''t believe in such a thing. So, my ghosts, relativity or math were back </w:t></w:r><w:r><w:rPr/><w:t>i</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">n the game. Something went wrong with my certificate from praxis, what was it? As I said, they suck. Finally, when my problem with the certificate were fix it, classes had started 3 days ago. </w:t></w:r><w:r><w:rPr/><w:t>And they DID NOT allow me to join the god damn group. I don''t know, maybe personal pronouns or the verb BE, was TOOO much to handle for me, ummm no, better, maybe their British accent will  be something impossible for me, and I will not have another option than cry because I couldn''t understand them, puffff. So, what ever the reason was, I couldn''t join them, and I didn''t have other option, than wait for the next group of brilliant, international and future British group. But, what I didn''t realize at that moment, is that in fact, that group would be special to me, and even craziest, there was a British soul.</w:t></w:r></w:p>'
$insPoint1.InsertXML($xml1)

# --- Edit 2: "That week, before classes started..." paragraph ---
# (re-locate: paragraph indices shifted after edit 1, but text offsets
# after it are unaffected in absolute terms only if we search by anchor.)
$target2 = Get-ParagraphByStart("That week, before classes started")
if ($null -eq $target2) {
    throw "Could not locate paragraph 2 (That week, before classes started...)"
}
$r2 = $target2.Range
$delRange2 = $d.Range($r2.Start, $r2.End + 1)
$delRange2.Delete()

$insPoint2 = $d.Range($r2.Start, $r2.Start)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">That week, before classes started, I was thinking about, why I had to do, whatever they wanna that i do. And considered the option to just drop out. But I didn''t. Why, the relativity, math, destiny or the British soul </w:t></w:r><w:r><w:rPr/><w:t>(British soul, is the name that I''m gonna use to refer to the British soul)</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">. I don''t know, I just know that I was there, waiting for “the big day”, for learn again about the verb BE and his friends, with people that just wanted to learn English, NOT, to love the English, which I think, it is the most important part of learning English </w:t></w:r><w:r><w:rPr/><w:t>or whatever you wanna learn</w:t></w:r><w:r><w:rPr/><w:t>. However, I had to go.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>So, June 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> of 2014. The day that I never gonna forget, why,  lets say, that, it was a mix of emotions, like a roller coaster, </w:t></w:r><w:r><w:rPr/><w:t>I was working on a company, that we are still trying to keep alive, by the way, and It was a long day, too much work  and best of all (*Irony injection), the past was calling me . God , you must be kidding me, I thought.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">6:00 pm. I watched the clock. (Classes were from 6:30 pm to 9:30 pm. ). </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">Sitting </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">At Udi''s cafeteria, next to my BF and telling him that I really did NOT want to go, I decided to go. But, I took my time to go. I went into the the </w:t></w:r><w:r><w:rPr/><w:t>class</w:t></w:r><w:r><w:rPr/><w:t>room, my first impression. The teacher, what a creepy guy. I didn''t even know what to think about him. Made me feel like, if we were in Halloween. I just thought, far from him, better I will. But then I remember that, we should</w:t></w:r><w:r><w:rPr/><w:t>n''t</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> judge a book by its cover. </w:t></w:r><w:r><w:rPr/><w:t>Really if we do this, judge books, we can change the whole history.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>The 2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> thing that I could notice was that, there were a lot of people, I mean, we were a big group. </w:t></w:r><w:r><w:rPr/><w:t>It looked like a 1</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> semester engineering class room; kids dreaming with create the next facebook. So, I thought: “brand new chusma, lets shine and have fun with this kids”. But this chusma, was kind of different. Why, we get there, I promise. For now, it was time for the personal presentation and I knew that I would enjoy it a LOT. Because there were two friends of mine, and lets say that English, hate them. And best of all, remember the  asshole who let me in  the middle of a storm. Well, that son of A BITCH was there, and I  was anxious to hear what the bastard was going to say, because  my soul, will laugh for a while. Indeed, it happened what I said. When purunga(Creepy teacher) ask him to introduce himself, the magic just came out. I''m actually sure that with every single “word” that he was telling, Mr. William Shakespeare wanted to return from the land of deaths and kill him. God, how happy I was in that moment.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">We listened to few of the new people (I was one of them, we were like 8) and others form the old class. And after that, listening the little Shakespeares, I thought “Scholarship, here we go again” </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">and I said: “It''s time to shine and show these boys, what''s English”. I did, I introduce myself  with my own style, which means that nobody could “handle me”, including purunga. Well, that''s what I thought at that moment. But, writing this, I''m sure that, there was someone who actually could “handle me”. Yes, the British soul, and I bet that in that moment, the British soul wanted to stand up and correct me. She always does. But,she didn''t. I don''t know why, maybe because the British soul was a new one too. However, I think it was a big mistake form  her. Because form that moment until the end of that course, I always was, one point up to her. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>Sitting there, thinking that I was in my kingdom, it was the time to realize the 3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> thing of that night, time to realize my British soul,ok ok, I hadn''t noticed yet, we were too many, remember?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>And she did. What she did, she defeated me, and I was playing local or that''s what I thought.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>I''m not good at remember things, specially with clothes. But, I can remember how she looked that day,  wearing a violet blouse and her typical blue jeans with an interesting boots. I liked, simple, but different  at the same time. But she looked pretty young, no id and like the others. So, I did NOT trust, even when I was feeling, that, she had something different, something that I would like to know. But the fact of no id, scared me  a lot. Don''t ask me why, lets say, ghosts problems.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>But then, she spoke and for a while, I have to say that I thought that I was in a dream or nightmare, because I saw my kingdom threatened. Actually, I saw many things threatened. And it was just the first day, the beginning of the class, the firsts words. My impulses just woke up, and started talking to me and I told them “shut the fuck up, and let me listen to her”. And that''s what I did, my eyes were deeper  on her with every single word, ''cause, damn, she did really well. And when she stopped talking, I had something very very clear in my mind and it was...</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>'
$insPoint2.InsertXML($xml2)

Write-Output "done"
